$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.770.41"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "1.614.60"
$ws.Range("E3").Value = "  -3.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.003"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.366"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08410"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.021"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001273"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "1.610.08"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06910"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.819"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("D24").Value = "23.780.28"
$ws.Range("E24").Value = "  -3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.860"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.233"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.80%  "
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.488"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "1.784.22"
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08020"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9728"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02874"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.553"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09190"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.417"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7466"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6854"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.458"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.054"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.203"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.02%  "
